$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.97750833333333
$ws.Range("H2").Value = 65.932525
$ws.Range("I2").Value = 0.5427578249542736
$ws.Range("J2").Value = 0.5427578249542736
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 40.314886
$ws.Range("N2").Value = 120.944658
$ws.Range("O2").Value = 0.1963336494301312
$ws.Range("P2").Value = 0.1963336494301312
$ws.Range("Q2").Value = 886.0207430223834
$ws.Range("R2").Value = 7974.18668720145
$ws.Range("S2").Value = 0.1065616245300328
$ws.Range("T2").Value = 0.1065616245300328
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.97750833333333
$ws.Range("H3").Value = 65.932525
$ws.Range("I3").Value = 0.5427578249542736
$ws.Range("J3").Value = 0.5427578249542736
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 69.65329366666667
$ws.Range("N3").Value = 208.959881
$ws.Range("O3").Value = 0.3392118072814421
$ws.Range("P3").Value = 0.3392118072814421
$ws.Range("Q3").Value = 1530.805842003281
$ws.Range("R3").Value = 13777.25257802952
$ws.Range("S3").Value = 0.1841098627188837
$ws.Range("T3").Value = 0.1841098627188837
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.97750833333333
$ws.Range("H4").Value = 65.932525
$ws.Range("I4").Value = 0.5427578249542736
$ws.Range("J4").Value = 0.5427578249542736
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 69.40355333333332
$ws.Range("N4").Value = 208.21066
$ws.Range("O4").Value = 0.3379955709003387
$ws.Range("P4").Value = 0.3379955709003388
$ws.Range("Q4").Value = 1525.317171746277
$ws.Range("R4").Value = 13727.8545457165
$ws.Range("S4").Value = 0.1834497409060458
$ws.Range("T4").Value = 0.1834497409060458
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.97750833333333
$ws.Range("H5").Value = 65.932525
$ws.Range("I5").Value = 0.5427578249542736
$ws.Range("J5").Value = 0.5427578249542736
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.96691433333334
$ws.Range("N5").Value = 77.90074300000001
$ws.Range("O5").Value = 0.126458972388088
$ws.Range("P5").Value = 0.126458972388088
$ws.Range("Q5").Value = 570.6880761517862
$ws.Range("R5").Value = 5136.192685366075
$ws.Range("S5").Value = 0.06863659679931118
$ws.Range("T5").Value = 0.06863659679931118
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.688376
$ws.Range("H6").Value = 38.065128
$ws.Range("I6").Value = 0.3133528721960219
$ws.Range("J6").Value = 0.3133528721960219
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.314886
$ws.Range("N6").Value = 120.944658
$ws.Range("O6").Value = 0.1963336494301312
$ws.Range("P6").Value = 0.1963336494301312
$ws.Range("Q6").Value = 511.530431965136
$ws.Range("R6").Value = 4603.773887686224
$ws.Range("S6").Value = 0.06152171295765846
$ws.Range("T6").Value = 0.06152171295765846
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.688376
$ws.Range("H7").Value = 38.065128
$ws.Range("I7").Value = 0.3133528721960219
$ws.Range("J7").Value = 0.3133528721960219
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 69.65329366666667
$ws.Range("N7").Value = 208.959881
$ws.Range("O7").Value = 0.3392118072814421
$ws.Range("P7").Value = 0.3392118072814421
$ws.Range("Q7").Value = 883.7871796810854
$ws.Range("R7").Value = 7954.084617129768
$ws.Range("S7").Value = 0.1062929940944433
$ws.Range("T7").Value = 0.1062929940944433
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.688376
$ws.Range("H8").Value = 38.065128
$ws.Range("I8").Value = 0.3133528721960219
$ws.Range("J8").Value = 0.3133528721960219
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 69.40355333333332
$ws.Range("N8").Value = 208.21066
$ws.Range("O8").Value = 0.3379955709003387
$ws.Range("P8").Value = 0.3379955709003388
$ws.Range("Q8").Value = 880.6183804293865
$ws.Range("R8").Value = 7925.565423864479
$ws.Range("S8").Value = 0.1059118829311553
$ws.Range("T8").Value = 0.1059118829311553
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.688376
$ws.Range("H9").Value = 38.065128
$ws.Range("I9").Value = 0.3133528721960219
$ws.Range("J9").Value = 0.3133528721960219
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.96691433333334
$ws.Range("N9").Value = 77.90074300000001
$ws.Range("O9").Value = 0.126458972388088
$ws.Range("P9").Value = 0.126458972388088
$ws.Range("Q9").Value = 329.4779726211227
$ws.Range("R9").Value = 2965.301753590104
$ws.Range("S9").Value = 0.0396262822127648
$ws.Range("T9").Value = 0.0396262822127648
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.008189666666667
$ws.Range("H10").Value = 6.024569
$ws.Range("I10").Value = 0.04959436889042158
$ws.Range("J10").Value = 0.04959436889042158
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 40.314886
$ws.Range("N10").Value = 120.944658
$ws.Range("O10").Value = 0.1963336494301312
$ws.Range("P10").Value = 0.1963336494301312
$ws.Range("Q10").Value = 80.95993747804467
$ws.Range("R10").Value = 728.639437302402
$ws.Range("S10").Value = 0.009737043435440632
$ws.Range("T10").Value = 0.009737043435440632
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.008189666666667
$ws.Range("H11").Value = 6.024569
$ws.Range("I11").Value = 0.04959436889042158
$ws.Range("J11").Value = 0.04959436889042158
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 69.65329366666667
$ws.Range("N11").Value = 208.959881
$ws.Range("O11").Value = 0.3392118072814421
$ws.Range("P11").Value = 0.3392118072814421
$ws.Range("Q11").Value = 139.8770245906988
$ws.Range("R11").Value = 1258.893221316289
$ws.Range("S11").Value = 0.01682299550230243
$ws.Range("T11").Value = 0.01682299550230243
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.008189666666667
$ws.Range("H12").Value = 6.024569
$ws.Range("I12").Value = 0.04959436889042158
$ws.Range("J12").Value = 0.04959436889042158
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 69.40355333333332
$ws.Range("N12").Value = 208.21066
$ws.Range("O12").Value = 0.3379955709003387
$ws.Range("P12").Value = 0.3379955709003388
$ws.Range("Q12").Value = 139.3754986339489
$ws.Range("R12").Value = 1254.37948770554
$ws.Range("S12").Value = 0.01676267702656004
$ws.Range("T12").Value = 0.01676267702656004
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.008189666666667
$ws.Range("H13").Value = 6.024569
$ws.Range("I13").Value = 0.04959436889042158
$ws.Range("J13").Value = 0.04959436889042158
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.96691433333334
$ws.Range("N13").Value = 77.90074300000001
$ws.Range("O13").Value = 0.126458972388088
$ws.Range("P13").Value = 0.126458972388088
$ws.Range("Q13").Value = 52.14648903941855
$ws.Range("R13").Value = 469.318401354767
$ws.Range("S13").Value = 0.006271652926118472
$ws.Range("T13").Value = 0.006271652926118472
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.818218
$ws.Range("H14").Value = 11.454654
$ws.Range("I14").Value = 0.09429493395928291
$ws.Range("J14").Value = 0.09429493395928291
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 40.314886
$ws.Range("N14").Value = 120.944658
$ws.Range("O14").Value = 0.1963336494301312
$ws.Range("P14").Value = 0.1963336494301312
$ws.Range("Q14").Value = 153.931023393148
$ws.Range("R14").Value = 1385.379210538332
$ws.Range("S14").Value = 0.01851326850699922
$ws.Range("T14").Value = 0.01851326850699922
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.818218
$ws.Range("H15").Value = 11.454654
$ws.Range("I15").Value = 0.09429493395928291
$ws.Range("J15").Value = 0.09429493395928291
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 69.65329366666667
$ws.Range("N15").Value = 208.959881
$ws.Range("O15").Value = 0.3392118072814421
$ws.Range("P15").Value = 0.3392118072814421
$ws.Range("Q15").Value = 265.9514596373527
$ws.Range("R15").Value = 2393.563136736174
$ws.Range("S15").Value = 0.03198595496581258
$ws.Range("T15").Value = 0.03198595496581258
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.818218
$ws.Range("H16").Value = 11.454654
$ws.Range("I16").Value = 0.09429493395928291
$ws.Range("J16").Value = 0.09429493395928291
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 69.40355333333332
$ws.Range("N16").Value = 208.21066
$ws.Range("O16").Value = 0.3379955709003387
$ws.Range("P16").Value = 0.3379955709003388
$ws.Range("Q16").Value = 264.9978966012933
$ws.Range("R16").Value = 2384.98106941164
$ws.Range("S16").Value = 0.03187127003657757
$ws.Range("T16").Value = 0.03187127003657757
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.818218
$ws.Range("H17").Value = 11.454654
$ws.Range("I17").Value = 0.09429493395928291
$ws.Range("J17").Value = 0.09429493395928291
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.96691433333334
$ws.Range("N17").Value = 77.90074300000001
$ws.Range("O17").Value = 0.126458972388088
$ws.Range("P17").Value = 0.126458972388088
$ws.Range("Q17").Value = 99.14733971199134
$ws.Range("R17").Value = 892.326057407922
$ws.Range("S17").Value = 0.01192444044989354
$ws.Range("T17").Value = 0.01192444044989354
